# Apply updated "want to go" counts (column F) across the relevant sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2166
$ws1.Range("F5").Value = 11430
$ws1.Range("F7").Value = 318
$ws1.Range("F9").Value = 11369
$ws1.Range("F10").Value = 462
$ws1.Range("F13").Value = 1747
$ws1.Range("F14").Value = 5667
$ws1.Range("F16").Value = 3484
$ws1.Range("F17").Value = 176

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2166
$ws4.Range("F7").Value = 11430
$ws4.Range("F9").Value = 318
$ws4.Range("F11").Value = 11369
$ws4.Range("F12").Value = 462
$ws4.Range("F15").Value = 1747
$ws4.Range("F16").Value = 1
$ws4.Range("F17").Value = 5667
$ws4.Range("F19").Value = 3484
$ws4.Range("F20").Value = 176
